$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 101. This pushes the existing
# rows 101-114 down to 103-116 (matching the new dimension A1:T116) and
# inherits the date number-format (style) on column D from the surrounding
# rows automatically.
$ws.Rows("101:102").Insert()

# New row 101: "Especial" quality entry for the 2023-04-06 (serial 45027) batch.
$ws.Range("A101").Value = 3
$ws.Range("B101").Value = "Femacal de La Calera"
$ws.Range("C101").Value = "Coquimbo"
$ws.Range("D101").Value = 45027
$ws.Range("E101").Value = 5
$ws.Range("F101").Value = "Fruta"
$ws.Range("G101").Value = 100104
$ws.Range("H101").Value = "Frutos de pepita"
$ws.Range("I101").Value = 100104003
$ws.Range("J101").Value = "Membrillo"
$ws.Range("K101").Value = "Champion"
$ws.Range("L101").Value = "Especial"
$ws.Range("M101").Value = 56
$ws.Range("N101").Value = 16000
$ws.Range("O101").Value = 16000
$ws.Range("P101").Value = 16000
$ws.Range("Q101").Value = "$/caja 18 kilos empedrada"
$ws.Range("R101").Value = "Región de O'Higgins"
$ws.Range("S101").Value = 889
$ws.Range("T101").Value = 18

# New row 102: "Primera" quality entry for the same 45027 batch.
$ws.Range("A102").Value = 3
$ws.Range("B102").Value = "Femacal de La Calera"
$ws.Range("C102").Value = "Coquimbo"
$ws.Range("D102").Value = 45027
$ws.Range("E102").Value = 5
$ws.Range("F102").Value = "Fruta"
$ws.Range("G102").Value = 100104
$ws.Range("H102").Value = "Frutos de pepita"
$ws.Range("I102").Value = 100104003
$ws.Range("J102").Value = "Membrillo"
$ws.Range("K102").Value = "Champion"
$ws.Range("L102").Value = "Primera"
$ws.Range("M102").Value = 60
$ws.Range("N102").Value = 14000
$ws.Range("O102").Value = 14000
$ws.Range("P102").Value = 14000
$ws.Range("Q102").Value = "$/caja 18 kilos empedrada"
$ws.Range("R102").Value = "Región de O'Higgins"
$ws.Range("S102").Value = 778
$ws.Range("T102").Value = 18
